# Generate Report for Handoff
# Update the "Overview" sheet status from "In Translation" to "Ready for handoff"
# and bump the handoff timestamps forward for the zh-cn locale row.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!B2 (zh-cn status) and Overview!C2 (de-de status)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# Overview!D2 - Latest Handoff Date
$overview.Range("D2").Value = "2016-03-21 18:36:16"

# de-de!E2 - Latest Handoff Datetime (same logical value as Overview!D2)
$dede.Range("E2").Value = "2016-03-21 18:36:16"

# zh-cn!E2 - Latest Handoff Datetime
$zhcn.Range("E2").Value = "2016-03-21 18:36:12"
